$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Save" column in H1, matching the formatting of the
# existing header cells (e.g. G1: bold font, border, centered alignment).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("H1").Value = "Save"

# Fill in the "Save" column values for rows 2-7
$values = @(1, 1, 0, 0, 1, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
